$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("other_resources")
$ws.Rows.Item(6).Insert()
$ws.Range("A6").Value = "Podcast"
$ws.Range("B6").Value = "TEMP"
try {
  $ws.Sort.SortFields.Clear()
  $ws.Sort.SortFields.Add($ws.Range("A2:A7"))
  $ws.Sort.SetRange($ws.Range("A2:B7"))
  $ws.Sort.Header = 0
  $ws.Sort.Apply()
  Write-Output "sort applied"
} catch {
  Write-Output "sort apply failed: $_"
}
$ws.Activate()
